$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B5 (marketDaysMode value) changes from "Manual" to "Auto"
$ws.Range("B5").Value = "Auto"
